$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "Jacopo Ricci"
$ws.Range("B27").Value = "Daniele  Dalbosco | iMontagna"
$ws.Range("C27").Value = "Federico Andreis | iMontagna"
$ws.Range("D27").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("E27").Value = "MARTINO TAMONI | U.S. Guarna"
$ws.Range("F27").Value = "Halzyd  Pupuleku | F.C. Sala Giardini"
